# predictions updated using weather data
#
# Adds 5 new rows (58-62) of weekly prediction data to Sheet1, following
# the same layout as the existing rows:
#   A: the day the prediction is made
#   B: the week of the target variable
#   C: Real
#   D: Prediction
#   E: difference
#   F: Model
#   G: MASE (test)
#   H: MAPE (test)
#   I: MAE (test)
#   J: Weekly MAE
#   K: Weekly MAPE

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$modelName = "KNN"

# Column A must hold the literal text "2021-01-09" (exactly like the
# existing rows above it), not an auto-converted date value. Typing that
# string directly into a cell makes Excel silently turn it into a date
# serial number, so instead we build the text with a TEXT() formula in a
# scratch area, copy it, and paste only the resulting values into the
# target cells; this keeps the cells as plain text using the default
# style, with no left-over number formatting.
$dateSerial = 44205  # serial number for 2021-01-09
$scratch = $ws.Range("Z1:Z5")
$scratch.Formula = "=TEXT(" + $dateSerial + ",""yyyy-mm-dd"")"
$scratch.Copy()
$ws.Range("A58:A62").PasteSpecial(-4163)  # xlPasteValues
$scratch.Clear()
$excel.CutCopyMode = $false

# Row 58 : full stats row for week "10 Jan -- 16 Jan 2021"
$ws.Range("B58").Value = "10 Jan -- 16 Jan 2021"
$ws.Range("C58").Value = 3333.57
$ws.Range("D58").Value = 1128.06
$ws.Range("E58").Value = 2205.51
$ws.Range("F58").Value = $modelName
$ws.Range("G58").Value = 0.87
$ws.Range("H58").Value = 59.66
$ws.Range("I58").Value = 1893.79
$ws.Range("J58").Value = 1761.8
$ws.Range("K58").Value = 56.92

# Row 59 : week "17 Jan -- 23 Jan 2021"
$ws.Range("B59").Value = "17 Jan -- 23 Jan 2021"
$ws.Range("D59").Value = 1266.5
$ws.Range("F59").Value = $modelName

# Row 60 : week "24 Jan -- 30 Jan 2021"
$ws.Range("B60").Value = "24 Jan -- 30 Jan 2021"
$ws.Range("D60").Value = 1274.17
$ws.Range("F60").Value = $modelName

# Row 61 : week "31 Jan -- 06 Feb 2021"
$ws.Range("B61").Value = "31 Jan -- 06 Feb 2021"
$ws.Range("D61").Value = 1286.63
$ws.Range("F61").Value = $modelName

# Row 62 : week "07 Feb -- 13 Feb 2021"
$ws.Range("B62").Value = "07 Feb -- 13 Feb 2021"
$ws.Range("D62").Value = 1526.29
$ws.Range("F62").Value = $modelName
